$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Legislature" row entirely (old row 12); everything below shifts up.
$ws.Rows("12:12").Delete()

# Groundwater Management (SWRCB) row (now row 20 after the shift):
# columns D and F change from -0.5 to -1
$ws.Range("D20").Value = -1
$ws.Range("F20").Value = -1

# Match the author's final on-screen selection (bottom-right pane of the
# frozen view) so the saved view state lines up with the source workbook.
$ws.Range("F21").Select()
